$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 390.6111
$ws.Range("I2").Value = 232.72728
$ws.Range("K2").Value = 232.72728
$ws.Range("M2").Value = -119.72728
$ws.Range("H4").Value = 660.2222
$ws.Range("I4").Value = 735
$ws.Range("J4").Value = 566.75
$ws.Range("K4").Value = 735
$ws.Range("L4").Value = 566.75
$ws.Range("M4").Value = -621
$ws.Range("N4").Value = -794.75
$ws.Range("H5").Value = 203.71428
$ws.Range("I5").Value = 71.166664
$ws.Range("J5").Value = 999
$ws.Range("K5").Value = 71.166664
$ws.Range("L5").Value = 999
$ws.Range("M5").Value = 43.833336
$ws.Range("N5").Value = -1229
$ws.Range("H18").Value = 816.375
$ws.Range("I18").Value = 816.375
$ws.Range("K18").Value = 816.375
$ws.Range("M18").Value = -532.375
$ws.Range("H62").Value = 8749.333000000001
$ws.Range("J62").Value = 8749.333000000001
$ws.Range("L62").Value = 8749.333000000001
$ws.Range("N62").Value = -9997.333000000001
$ws.Range("H65").Value = 8749.333000000001
$ws.Range("J65").Value = 8749.333000000001
$ws.Range("L65").Value = 43746.665
$ws.Range("N65").Value = -49986.665
$ws.Range("H87").Value = 58569.168
$ws.Range("J87").Value = 67883
$ws.Range("L87").Value = 67883
$ws.Range("N87").Value = -70379
$ws.Range("H90").Value = 58569.168
$ws.Range("J90").Value = 67883
$ws.Range("L90").Value = 203649
$ws.Range("N90").Value = -216129
$ws.Range("H100").Value = 5832.8887
$ws.Range("I100").Value = 3874
$ws.Range("K100").Value = 3874
$ws.Range("M100").Value = -3333

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H101").Value = 71131
$ws.Range("J101").Value = 71131
$ws.Range("L101").Value = 71131
$ws.Range("N101").Value = -77621
$ws.Range("H103").Value = 2362
$ws.Range("J103").Value = 2362
$ws.Range("L103").Value = 2362
$ws.Range("N103").Value = -4706
$ws.Range("H106").Value = 0
$ws.Range("I106").Value = 0
$ws.Range("K106").Value = 0
$ws.Range("M106").ClearContents()
$ws.Range("H132").Value = 3343.375
$ws.Range("I132").Value = 2678.1904
$ws.Range("J132").Value = 7999.6665
$ws.Range("K132").Value = 8034.5712
$ws.Range("L132").Value = 23998.9995
$ws.Range("M132").Value = -5504.5712
$ws.Range("N132").Value = -29058.9995

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H16").Value = 1137.5
$ws.Range("I16").Value = 1137.5
$ws.Range("K16").Value = 1137.5
$ws.Range("M16").Value = -967.5
$ws.Range("H22").Value = 296.35715
$ws.Range("I22").Value = 293.91666
$ws.Range("K22").Value = 293.91666
$ws.Range("M22").Value = -120.91666
$ws.Range("H107").Value = 2669.389
$ws.Range("I107").Value = 945.2917
$ws.Range("K107").Value = 945.2917
$ws.Range("M107").Value = 974.7083
$ws.Range("H114").Value = 150000
$ws.Range("J114").Value = 150000
$ws.Range("L114").Value = 150000
$ws.Range("N114").Value = -158678
$ws.Range("H134").Value = 5811.8667
$ws.Range("I134").Value = 5811.8667
$ws.Range("K134").Value = 17435.6001
$ws.Range("M134").Value = -14900.6001
$ws.Range("H140").Value = 123259.664
$ws.Range("I140").Value = 99999
$ws.Range("K140").Value = 99999
$ws.Range("M140").Value = -94819

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 85.2
$ws.Range("I7").Value = 85.818184
$ws.Range("K7").Value = 85.818184
$ws.Range("M7").Value = 27.181816
$ws.Range("H12").Value = 647.1111
$ws.Range("I12").Value = 238.16667
$ws.Range("J12").Value = 1465
$ws.Range("K12").Value = 238.16667
$ws.Range("L12").Value = 1465
$ws.Range("M12").Value = -68.16667000000001
$ws.Range("N12").Value = -1805
$ws.Range("H52").Value = 268975.56
$ws.Range("I52").Value = 40000
$ws.Range("J52").Value = 334397.16
$ws.Range("K52").Value = 40000
$ws.Range("L52").Value = 334397.16
$ws.Range("M52").Value = -39706
$ws.Range("N52").Value = -334985.16
$ws.Range("H81").Value = 46246.75
$ws.Range("J81").Value = 46662.332
$ws.Range("L81").Value = 46662.332
$ws.Range("N81").Value = -48658.332
$ws.Range("H84").Value = 46246.75
$ws.Range("J84").Value = 46662.332
$ws.Range("L84").Value = 139986.996
$ws.Range("N84").Value = -149970.996
$ws.Range("H94").Value = 5034.25
$ws.Range("I94").Value = 2065.25
$ws.Range("K94").Value = 2065.25
$ws.Range("M94").Value = -1614.25
$ws.Range("H122").Value = 1739
$ws.Range("I122").Value = 1691.6875
$ws.Range("K122").Value = 5075.0625
$ws.Range("M122").Value = -2625.0625

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1871.1428
$ws.Range("I5").Value = 772
$ws.Range("K5").Value = 2316
$ws.Range("M5").Value = -2204
$ws.Range("H12").Value = 147
$ws.Range("I12").Value = 114.6
$ws.Range("J12").Value = 163.2
$ws.Range("K12").Value = 343.8
$ws.Range("L12").Value = 489.6
$ws.Range("M12").Value = -170.8
$ws.Range("N12").Value = -835.5999999999999
$ws.Range("H14").Value = 3500.6667
$ws.Range("I14").Value = 3500.6667
$ws.Range("K14").Value = 10502.0001
$ws.Range("M14").Value = -10329.0001
$ws.Range("H40").Value = 79.181816
$ws.Range("I40").Value = 49.5
$ws.Range("J40").Value = 114.8
$ws.Range("K40").Value = 198
$ws.Range("L40").Value = 459.2
$ws.Range("M40").Value = -129
$ws.Range("N40").Value = -597.2
$ws.Range("H104").Value = 6000
$ws.Range("I104").Value = 4666.6665
$ws.Range("K104").Value = 13999.9995
$ws.Range("M104").Value = -11378.9995
$ws.Range("H108").Value = 1650.25
$ws.Range("I108").Value = 926
$ws.Range("K108").Value = 2778
$ws.Range("M108").Value = 102
$ws.Range("H109").Value = 145229.42
$ws.Range("J109").Value = 3277.5
$ws.Range("L109").Value = 9832.5
$ws.Range("N109").Value = -11912.5
$ws.Range("H135").Value = 1871.1428
$ws.Range("I135").Value = 772
$ws.Range("K135").Value = 6948
$ws.Range("M135").Value = -4413

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 112.35714
$ws.Range("J2").Value = 139.55556
$ws.Range("L2").Value = 139.55556
$ws.Range("N2").Value = -365.55556
$ws.Range("H100").Value = 33000
$ws.Range("J100").Value = 33000
$ws.Range("L100").Value = 33000
$ws.Range("N100").Value = -35164
$ws.Range("H135").Value = 235000
$ws.Range("J135").Value = 235000
$ws.Range("L135").Value = 235000
$ws.Range("N135").Value = -245140
$ws.Range("H138").Value = 75000
$ws.Range("I138").Value = 75000
$ws.Range("K138").Value = 75000
$ws.Range("M138").Value = -69860

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H108").Value = 0
$ws.Range("J108").Value = 0
$ws.Range("L108").Value = 0
$ws.Range("N108").ClearContents()
$ws.Range("H135").Value = 88375
$ws.Range("J135").Value = 88375
$ws.Range("L135").Value = 88375
$ws.Range("N135").Value = -98515

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 49662.668
$ws.Range("J46").Value = 49662.668
$ws.Range("L46").Value = 49662.668
$ws.Range("N46").Value = -50124.668
$ws.Range("H94").Value = 5324999
$ws.Range("J94").Value = 5324999
$ws.Range("L94").Value = 5324999
$ws.Range("N94").Value = -5326801
$ws.Range("H103").Value = 40511.2
$ws.Range("J103").Value = 40511.2
$ws.Range("L103").Value = 40511.2
$ws.Range("N103").Value = -42855.2
$ws.Range("H107").Value = 733.5
$ws.Range("I107").Value = 543.0769
$ws.Range("K107").Value = 1629.2307
$ws.Range("M107").Value = 290.7692999999999
$ws.Range("H126").Value = 3753.3157
$ws.Range("I126").Value = 2339.5386
$ws.Range("J126").Value = 6816.5
$ws.Range("K126").Value = 7018.6158
$ws.Range("L126").Value = 20449.5
$ws.Range("M126").Value = -4548.6158
$ws.Range("N126").Value = -25389.5
$ws.Range("H134").Value = 49662.668
$ws.Range("J134").Value = 49662.668
$ws.Range("L134").Value = 148988.004
$ws.Range("N134").Value = -154058.004
$ws.Range("H135").Value = 35941.8
$ws.Range("J135").Value = 35941.8
$ws.Range("L135").Value = 35941.8
$ws.Range("N135").Value = -46081.8
